# Advanced/一致性哈希.docx edit
#
# The commit fixes a word-order typo in the "虚拟节点" (virtual node)
# paragraph: "一实际个节点" -> "一个实际节点", and (as a side effect of
# Word's edit) the "实际" portion becomes its own run and the document's
# lone "_GoBack" bookmark (which previously sat near "甚至更大，" at the
# very end of the document) moves to sit right after that "实际" run.

$d = $word.ActiveDocument

# --- 1. Fix the word order: "一实际个节点" -> "一个实际节点" -------------
$rng = $d.Content
$rng.Find.ClearFormatting()
$found = $rng.Find.Execute("一实际个节点", $true, $false, $false, $false, $false, $true, 1, $false, "一个实际节点", 2)

# --- 2. Re-locate the fixed text and split "实际" into its own run, -----
#        moving the (only) "_GoBack" bookmark to sit right after it.
$rng2 = $d.Content
$rng2.Find.ClearFormatting()
$rng2.Find.Execute("一个实际节点") | Out-Null
$segStart = $rng2.Start

# offsets within "一个实际节点": 一(0) 个(1) 实(2) 际(3) 节(4) 点(5)
$shiShiStart = $segStart + 2
$shiShiEnd = $segStart + 4

# Touch-then-revert a character property on just "实际" so the engine
# splits it into its own <w:r> (formatting ends up identical to its
# neighbours, matching the target XML).
$shiji = $d.Range($shiShiStart, $shiShiEnd)
$shiji.Font.Bold = $true
$shiji.Font.Bold = $false

# Move the document's "_GoBack" bookmark to the boundary right after
# "实际" (Bookmarks.Add with an existing name re-seats that bookmark,
# which also removes it from its old location near "甚至更大，").
$bookmarkSpot = $d.Range($shiShiEnd, $shiShiEnd)
$d.Bookmarks.Add("_GoBack", $bookmarkSpot) | Out-Null
